{"js": "// The project was renamed from \"Code Royale\" to \"GroofyCode\" throughout the\n// proposal, and the Introduction paragraph's \"frontend\" was tightened to the\n// hyphenated \"front-end\". Apply both as scoped, case-sensitive find/replace\n// passes so every other run (formatting, list items, etc.) is left intact.\n\nconst body = context.document.body;\n\nconst nameHits = body.search(\"Code Royale\", { matchCase: true, matchWholeWord: false });\nnameHits.load(\"items\");\nconst frontendHits = body.search(\"frontend\", { matchCase: true, matchWholeWord: false });\nfrontendHits.load(\"items\");\n\nawait context.sync();\n\nfor (let i = 0; i < nameHits.items.length; i++) {\n  nameHits.items[i].insertText(\"GroofyCode\", \"Replace\");\n}\n\nfor (let i = 0; i < frontendHits.items.length; i++) {\n  frontendHits.items[i].insertText(\"front-end\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The project was renamed from \"Code Royale\" to \"GroofyCode\" throughout the\n# proposal, and the Introduction paragraph's \"frontend\" was tightened to the\n# hyphenated \"front-end\". Use Find/Replace (ReplaceAll) scoped to the whole\n# document body so every other run (bold headings, bullet lists, etc.) is\n# left untouched.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Code Royale\",  # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap -> wdFindContinue\n    $false,         # Format\n    \"GroofyCode\",   # ReplaceWith\n    2               # Replace -> wdReplaceAll\n)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"frontend\",     # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap -> wdFindContinue\n    $false,         # Format\n    \"front-end\",    # ReplaceWith\n    2               # Replace -> wdReplaceAll\n)\n"}
